$d = $word.ActiveDocument

# 1) Trim the sentence that used to run on after "-0.4120. " down to just that lead-in.
$rsq = [char]0x2019
$oldText = "0.02932, it tells me the root is -0.4120. It appears that these cases, where it doesn" + $rsq + "t find the closest root, happen when we are almost in the middle of two roots. "
$newText = "0.02932, it tells me the root is -0.4120. "
$found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# 2) Insert the new explanatory text right before the _GoBack bookmark (which sits
#    immediately after the just-shortened sentence), building it up piece by piece so the
#    apostrophe in "y'" can later be marked as its own superscript run.
$bm = $d.Bookmarks("_GoBack")
$insPos = $bm.Start
$ins = $d.Range($insPos, $insPos)
$ins.InsertAfter("This happens when we are near a maximum or minimum value of the function. In the case above, it skips to the -0.4120 value because it is a point near a maximum, so the slope is not very steep. So, when it does the necessary calculations to find the next X value, because y")

$bm = $d.Bookmarks("_GoBack")
$insPos = $bm.Start
$supStart = $insPos
$supRange = $d.Range($insPos, $insPos)
$supRange.InsertAfter($rsq)

$bm = $d.Bookmarks("_GoBack")
$insPos = $bm.Start
$tailRange = $d.Range($insPos, $insPos)
$tailRange.InsertAfter(" is closer to 1, it goes back farther and it skips over the closest root it had previously been attaching to. ")

# Now that all the surrounding text is in place, mark just the apostrophe as superscript.
$supRange = $d.Range($supStart, $supStart + 1)
$supRange.Font.Superscript = $true

# 3) Drop the old trailing single-space run that used to follow the bookmark.
$bm = $d.Bookmarks("_GoBack")
$afterBm = $bm.End
$tail = $d.Range($afterBm, $afterBm + 1)
if ($tail.Text -eq " ") {
    $tail.Delete()
}
